$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 3
$ws.Range("B8").Value = 4

$ws.Range("B5").Select()
